# Deletes the "photograph" row (row 133) from the wordlist so that the
# remaining rows (134:146) shift up by one, matching the reordering
# described in #260.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(133).Delete()
